$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The debtor table (rows 15-17) lists one row per overdue period (2507, 2508).
# A new overdue period (2509) needs to be added as an additional row, right
# below the existing rows, before the signature/footer block further down.
# Insert a new row 18 (this pushes the footer rows 22/23 down to 23/24).
$ws.Rows("18").Insert()

# Clone row 17 (which currently carries the "bottom of table" bordered style)
# down into the newly inserted row 18, together with its values - this will
# become the new last row of the table.
$ws.Range("B17:J17").Copy($ws.Range("B18:J18"))

# Row 17 is no longer the last row of the table, so give it the same
# "middle of table" style that row 16 already uses.
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)

# Row 18 represents the new period, 2509.
$ws.Range("E18").Value = "2509"

# Update the totals: the overdue amount (Valor Mora) and the period count
# (Cant. Periodos) both grew because of the new record above.
$ws.Range("E11").Value = 170820
$ws.Range("F13").Value = 3
